$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This worksheet is a syllabus table. The edit fixes a set of rows (13-23)
# whose B/C (value) cells were misaligned relative to their A (label) cells
# - inserting a dedicated data row for "Docentes responsaveis:" (the teacher
# name had been crammed into the "Objetivos:" row), filling in several
# previously-empty Portuguese fields (Objetivos, Programa resumido, Programa,
# Bibliografia) and appending a final row for the "Requisitos:" data.
# ---------------------------------------------------------------------------

# --- 1. Update cells that already exist (same address before/after) -------
$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "Fornecer ao estudante noções básicas de dispositivos digitais e suas aplicações com ênfase em microcontroladores e processadores digitais de sinais."
$ws.Range("C10").Value = "Fornecer ao estudante noções básicas de dispositivos digitais e suas aplicações com ênfase em microcontroladores e processadores digitais de sinais."

$ws.Range("A11").Value = "Objectives:"
$ws.Range("B11").Value = "Provide the student with the basics of digital devices and their applications with an emphasis on microcontrollers and digital signal processors."
$ws.Range("C11").Value = "Provide the student with the basics of digital devices and their applications with an emphasis on microcontrollers and digital signal processors."

$ws.Range("A12").Value = "Docentes responsáveis:"

$ws.Range("B13").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C13").Value = "519033 - Carlos Yujiro Shigue"

$ws.Range("A14").Value = "Programa resumido:"
$ws.Range("B14").Value = "Circuitos digitais. Microprocessadores e microcontroladores. Programação de sistemas de aquisição de dados e algoritmos de controle."
$ws.Range("C14").Value = "Circuitos digitais. Microprocessadores e microcontroladores. Programação de sistemas de aquisição de dados e algoritmos de controle."

$ws.Range("A15").Value = "Short syllabus:"
$ws.Range("B15").Value = "Digital circuits. Microprocessors and microcontrollers. Programming of data acquisition systems and control algorithms."
$ws.Range("C15").Value = "Digital circuits. Microprocessors and microcontrollers. Programming of data acquisition systems and control algorithms."

$ws.Range("A16").Value = "Programa:"
$ws.Range("B16").Value = "Bases numéricas. Aritmética binária. Funções lógicas. Álgebra de Boole. Minimização. Circuitos combinatórios. Flip-flops. Contadores e projeto de contadores. Introdução aos circuitos sequenciais. Microprocessadores. Microcontroladores e sistemas embarcados. Interfaces de comunicação. Linguagem de programação de baixo e alto nível na computação em tempo real. Desenvolvimento de protocolos de comando digital. Projeto com dispositivos programáveis: microcontroladores e processadores de sinais digitais. Programação de dispositivos FPGA."
$ws.Range("C16").Value = "Bases numéricas. Aritmética binária. Funções lógicas. Álgebra de Boole. Minimização. Circuitos combinatórios. Flip-flops. Contadores e projeto de contadores. Introdução aos circuitos sequenciais. Microprocessadores. Microcontroladores e sistemas embarcados. Interfaces de comunicação. Linguagem de programação de baixo e alto nível na computação em tempo real. Desenvolvimento de protocolos de comando digital. Projeto com dispositivos programáveis: microcontroladores e processadores de sinais digitais. Programação de dispositivos FPGA."

$ws.Range("A17").Value = "Syllabus:"

$ws.Range("A18").Value = "Avaliação:"

$ws.Range("A19").Value = "Método:"
$ws.Range("B19").Value = "Aulas expositivas, exercícios em sala, lista de exercícios, utilização de um simulador de circuitos, projeto de circuitos e atividades práticas em laboratório."
$ws.Range("C19").Value = "Aulas expositivas, exercícios em sala, lista de exercícios, utilização de um simulador de circuitos, projeto de circuitos e atividades práticas em laboratório."

$ws.Range("A20").Value = "Critério:"
$ws.Range("B20").Value = "Média ponderada de duas provas escritas, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + 2P2 + TR)/4"
$ws.Range("C20").Value = "Média ponderada de duas provas escritas, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + 2P2 + TR)/4"

$ws.Range("A21").Value = "Norma de recuperação:"
$ws.Range("B21").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Range("C21").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"

$ws.Range("A22").Value = "Bibliografia:"

# --- 2. Remove cells that no longer exist at this address ------------------
$ws.Range("A13").Clear()
$ws.Range("B18").Clear()
$ws.Range("C18").Clear()
$ws.Range("B23").Clear()
$ws.Range("C23").Clear()

# --- 3. Add brand-new cells (copy formatting from an already-correctly-
#        styled neighbour so no new style entries get minted) --------------
$ws.Range("B17").Value = "Numerical bases. Binary arithmetic. Logical functions. Boolean algebra. Minimization. Combinatorial circuits. flip-flops. Accountants and Accountants Design. Introduction to sequential circuits. Microprocessors. Microcontrollers and embedded systems. Communication interfaces. Low-level and high-level programming language in real-time computing. Development of digital command protocols. Project with programmable devices: microcontrollers and digital signal processors. Programming of FPGA devices."
$ws.Range("B16").Copy()
$ws.Range("B17").PasteSpecial(-4122)

$ws.Range("C17").Value = "Numerical bases. Binary arithmetic. Logical functions. Boolean algebra. Minimization. Combinatorial circuits. flip-flops. Accountants and Accountants Design. Introduction to sequential circuits. Microprocessors. Microcontrollers and embedded systems. Communication interfaces. Low-level and high-level programming language in real-time computing. Development of digital command protocols. Project with programmable devices: microcontrollers and digital signal processors. Programming of FPGA devices."
$ws.Range("C16").Copy()
$ws.Range("C17").PasteSpecial(-4122)

$ws.Range("B22").Value = "GAJSKI, D. D. Principles of Digital Design, Prentice Hall, 1997.`nTAUB, H. Circuitos Digitais e Microprocessadores, McGraw Hill, 1984.`nTOCCI, R. J.; AMBROSIO, F. J. Microprocessors and Microcomputers: Hardware and Software, Prentice Hall, 2002.`nCATSOULIS, J. Designing Embedded Hardware, OReilly Media, 2005.`nCRISP, J. Introduction to Microprocessors, Newnes, 2004.`nWILMSHURST, T. Designing Embedded Systems with PIC Microcontrollers, Newnes, 2009.`nDUBEY, R. Introduction to Embedded System Design using Field Programmable Gate Arrays, Springer, 2008.`nBATEMAN, A.; PATERSON-STEPHENS, I. The DSP Handbook: Algorithms, Applications and Design Techniques, Prentice Hall, 2002."
$ws.Range("B21").Copy()
$ws.Range("B22").PasteSpecial(-4122)

$ws.Range("C22").Value = "GAJSKI, D. D. Principles of Digital Design, Prentice Hall, 1997.`nTAUB, H. Circuitos Digitais e Microprocessadores, McGraw Hill, 1984.`nTOCCI, R. J.; AMBROSIO, F. J. Microprocessors and Microcomputers: Hardware and Software, Prentice Hall, 2002.`nCATSOULIS, J. Designing Embedded Hardware, OReilly Media, 2005.`nCRISP, J. Introduction to Microprocessors, Newnes, 2004.`nWILMSHURST, T. Designing Embedded Systems with PIC Microcontrollers, Newnes, 2009.`nDUBEY, R. Introduction to Embedded System Design using Field Programmable Gate Arrays, Springer, 2008.`nBATEMAN, A.; PATERSON-STEPHENS, I. The DSP Handbook: Algorithms, Applications and Design Techniques, Prentice Hall, 2002."
$ws.Range("C21").Copy()
$ws.Range("C22").PasteSpecial(-4122)

$ws.Range("A23").Value = "Requisitos:"
$ws.Range("A22").Copy()
$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("A23").Value = "Requisitos:"

$ws.Range("B24").Value = "LOM3263 -  Eletrônica Fundamental e Aplicada  (Requisito)`n"
$ws.Range("B19").Copy()
$ws.Range("B24").PasteSpecial(-4122)
$ws.Range("B24").Value = "LOM3263 -  Eletrônica Fundamental e Aplicada  (Requisito)`n"

$ws.Range("C24").Value = "LOM3263 -  Eletrônica Fundamental e Aplicada  (Requisito)`n"
$ws.Range("C19").Copy()
$ws.Range("C24").PasteSpecial(-4122)
$ws.Range("C24").Value = "LOM3263 -  Eletrônica Fundamental e Aplicada  (Requisito)`n"

# --- 4. Row heights ---------------------------------------------------------
$ws.Rows.Item(13).AutoFit()
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(18).AutoFit()
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 120
$ws.Rows.Item(23).AutoFit()
$ws.Rows.Item(24).RowHeight = 30
